{"js": "// \"Changes in the report scope\"\n// 1) Remove the yellow highlight from the \"Scope\" heading run.\n// 2) Rewrite the scope description sentence.\n// 3) Move the \"_GoBack\" bookmark from after \"Internet connectivity\" to just\n//    before \"are:\" inside the rewritten scope sentence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet scopeHeadingPara = null;\nlet scopeSentencePara = null;\nlet internetPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === \"Scope\") {\n    scopeHeadingPara = paragraphs.items[i];\n  }\n  if (text.indexOf(\"implement an optimized Question Answer matching system\") !== -1) {\n    scopeSentencePara = paragraphs.items[i];\n  }\n  if (text.indexOf(\"Internet connectivity\") !== -1) {\n    internetPara = paragraphs.items[i];\n  }\n}\n\nif (!scopeHeadingPara) throw new Error(\"Could not find the 'Scope' heading paragraph.\");\nif (!scopeSentencePara) throw new Error(\"Could not find the scope description paragraph.\");\nif (!internetPara) throw new Error(\"Could not find the 'Internet connectivity' paragraph.\");\n\n// 1) Strip the yellow highlight from \"Scope\".\nscopeHeadingPara.getRange().font.highlightColor = null;\n\n// 2) Rewrite the scope sentence.\nconst newSentence =\n  \"To implement an optimize Question Answer matching system having a better \" +\n  \"accuracy and a higher degree of similarity. The aspects to be under \" +\n  \"consideration are:\";\nscopeSentencePara.getRange().insertText(newSentence, \"Replace\");\nawait context.sync();\n\n// 3) Move the _GoBack bookmark into the rewritten sentence, right before \"are:\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst beforeAreColon = body.search(\"to be under consideration \", { matchCase: false });\nbeforeAreColon.load(\"items\");\nawait context.sync();\n\nif (beforeAreColon.items.length === 0) {\n  throw new Error(\"Could not find the bookmark insertion point.\");\n}\n\nconst bookmarkSpot = beforeAreColon.items[0].getRange(\"End\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"Changes in the report scope\"\n# 1) Remove the yellow highlight from the \"Scope\" heading run.\n# 2) Rewrite the scope description sentence.\n# 3) Move the \"_GoBack\" bookmark from after \"Internet connectivity\" to just\n#    before \"are:\" inside the rewritten scope sentence.\n\n$d = $word.ActiveDocument\n\n$newSentence = \"To implement an optimize Question Answer matching system having a better accuracy and a higher degree of similarity. The aspects to be under consideration are:\"\n\n# 1) Strip the yellow highlight from the \"Scope\" heading paragraph.\n$scopeFound = $false\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Scope\") {\n        $p.Range.Font.HighlightColorIndex = \"wdNoHighlight\"\n        $scopeFound = $true\n    }\n}\nif (-not $scopeFound) {\n    throw \"Could not find the 'Scope' heading paragraph.\"\n}\n\n# 2) Rewrite the scope sentence paragraph (replace text, keep paragraph mark).\n$sentenceFound = $false\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*implement an optimized*\") {\n        $r = $p.Range\n        [void]$r.MoveEnd(1, -1)\n        $r.Text = $newSentence\n        $sentenceFound = $true\n    }\n}\nif (-not $sentenceFound) {\n    throw \"Could not find the scope description paragraph.\"\n}\n\n# 3) Move the _GoBack bookmark into the rewritten sentence, right before \"are:\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$located = $findRange.Find.Execute(\"to be under consideration \")\nif (-not $located) {\n    throw \"Could not find the bookmark insertion point.\"\n}\n\n$bmRange = $findRange.Duplicate\n$bmRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
